$d = $word.ActiveDocument

# --- 1. Paragraph 5: merge the three runs (split apart by a gramStart/gramEnd
#        proofErr pair around "take into account") back into one plain run. ---
$p5 = $d.Paragraphs.Item(5)
$p5Start = $p5.Range.Start
$p5End = $p5.Range.End
$p5Range = $d.Range($p5Start, $p5End)
$p5Range.Text = "I defined access to public transit as a person living within a 10-minute walk of any bus stop. This is a very liberal definition" + [char]0x2014 + "many analyses use a 5-minute buffer around stops. This also does not take into account the requirement of transfers, which greatly affect the ease of use for many people outside the urban core."

# --- 2. Paragraph 7 ("Going forward ...") becomes the new "second analysis" paragraph. ---
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Text = "The second analysis looks at the workers who most likely still need to utilize public transit" + [char]0x2014 + "healthcare and service workers. Workers in these industries represent the few that still need to get to work during the shutdown, meaning the Port Authority should ensure that services cuts do not overly impact those workers."

# --- 3. Insert two brand-new paragraphs right after paragraph 7. ---
$insertOffset = $d.Paragraphs.Item(7).Range.End
$insertPoint = $d.Range($insertOffset, $insertOffset)

$apos = [char]0x2019
$emdash = [char]0x2014

$run1 = "This analysis centers around determining where healthcare and service workers respectively are relatively abundant. I" + $apos + "ve done this by constructing a relative abundance index that takes into account both the density of the type of worker in a give census block group, and the proportion of people in that block group that work in a given sector. In this way the relative abundance measure incorporates "
$run2 = "not only how many workers reside in a given block group, but also highlights where those workers are more likely to live."
$para9 = "Third, this project takes in bus usage data from the time of the bus changes, 3-25-20, until 4-21-20, and compares that data to the change in bus frequency. For this analysis I" + $apos + "ve constructed a ratio that compares that change in use to the change in buses" + $emdash + "values above 1 means use has decreased by less than the bus availability has decreased, or that use has actually increased. Values below 1 represent where buses were reduced more than use declined, and these values should receive scrutiny."

function XmlEscape([string]$s) {
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

# Only runs whose text has leading/trailing whitespace need xml:space="preserve"
# (matches how Word itself marks runs in the target revision).
function RunXml([string]$s) {
    $needsPreserve = ($s.Length -gt 0) -and (($s[0] -eq ' ') -or ($s[$s.Length - 1] -eq ' '))
    if ($needsPreserve) {
        return '<w:r><w:t xml:space="preserve">' + (XmlEscape $s) + '</w:t></w:r>'
    } else {
        return '<w:r><w:t>' + (XmlEscape $s) + '</w:t></w:r>'
    }
}

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    '<w:p>' + (RunXml $run1) + (RunXml $run2) + '</w:p>' +
    '<w:p>' + (RunXml $para9) + '</w:p>' +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($xml) | Out-Null

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
